$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the "last updated" date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "FPIEBP" sheet: update balancing priorities for "hard coal" (row 3) ---
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# Keep FPIEBP as the active/selected sheet, move the selection to E3
$wsFpiebp.Activate() | Out-Null
$wsFpiebp.Range("E3").Select() | Out-Null
